# Apply the updated cryptocurrency Price (D) and Volume(1h) (E) figures.
# Price cells whose new text would otherwise be auto-parsed as a number
# (plain decimals such as "239.76") are forced to keep a Text format so
# the literal string survives (matching values such as "1.10" or "0.0944").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.932.12"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "2.204.94"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.76"
$ws.Range("E5").Value = "  -2.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.623"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.67"
$ws.Range("E7").Value = "  -2.12%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.61"
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0944"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.04"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").Value = "2.533.39"
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.12"
$ws.Range("E15").Value = "  -2.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.835"
$ws.Range("E16").Value = "  -2.15%  "
$ws.Range("D17").Value = "2.188.43"
$ws.Range("E17").Value = "  -3.43%  "
$ws.Range("D18").Value = "41.792.01"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000107"
$ws.Range("E19").Value = "  +8.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.35"
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.09"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("E22").Value = "  +19.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "228.63"
$ws.Range("E23").Value = "  -0.72%  "
$ws.Range("E24").Value = "  -7.69%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.44"
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.24"
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("E29").Value = "  +3.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.86"
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.46"
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.53"
$ws.Range("E32").Value = "  +7.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0787"
$ws.Range("E33").Value = "  -3.99%  "
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("E35").Value = "  -5.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.111"
$ws.Range("E36").Value = "  -7.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.24"
$ws.Range("E37").Value = "  -5.86%  "
$ws.Range("E38").Value = "  -2.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.88"
$ws.Range("E39").Value = "  -6.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "65.03"
$ws.Range("E40").Value = "  +4.64%  "
$ws.Range("E41").Value = "  -3.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.59"
$ws.Range("E42").Value = "  -3.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.197"
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("E44").Value = "  +0.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.43"
$ws.Range("E45").Value = "  -3.88%  "
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.37"
$ws.Range("E47").Value = "  +4.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.10"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "2.410.51"
$ws.Range("E51").Value = "  -1.68%  "
